$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial value that was updated
# for every data row (rows 2 through 200) from 45172 (2023-09-03)
# to 45175 (2023-09-06).
$ws.Range("C2:C200").Value = 45175
